$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 97
$ws.Range("E3").Value = 39
$ws.Range("E10").Value = 552
$ws.Range("F10").Value = 262
$ws.Range("H10").Value = 357
$ws.Range("E11").Value = 354
$ws.Range("F11").Value = 189
$ws.Range("G11").Value = 65
$ws.Range("H11").Value = 254
$ws.Range("E12").Value = 544
$ws.Range("F12").Value = 285
$ws.Range("G12").Value = 86
$ws.Range("H12").Value = 371
$ws.Range("E16").Value = 203
$ws.Range("F16").Value = 103
$ws.Range("H16").Value = 151
$ws.Range("E22").Value = 169
$ws.Range("E24").Value = 210
$ws.Range("E25").Value = 266
$ws.Range("F25").Value = 129
$ws.Range("H25").Value = 189
$ws.Range("E26").Value = 156
$ws.Range("F26").Value = 95
$ws.Range("H26").Value = 120
$ws.Range("E27").Value = 328
$ws.Range("E28").Value = 197
$ws.Range("E29").Value = 168
$ws.Range("F29").Value = 94
$ws.Range("H29").Value = 135
$ws.Range("E30").Value = 208
$ws.Range("F30").Value = 121
$ws.Range("H30").Value = 173
$ws.Range("E32").Value = 182
$ws.Range("F32").Value = 108
$ws.Range("H32").Value = 146
$ws.Range("E33").Value = 293
$ws.Range("E34").Value = 216
$ws.Range("F34").Value = 142
$ws.Range("H34").Value = 181
$ws.Range("E35").Value = 149
$ws.Range("F35").Value = 92
$ws.Range("H35").Value = 119
$ws.Range("E36").Value = 72
$ws.Range("F36").Value = 40
$ws.Range("H36").Value = 50
$ws.Range("E37").Value = 159
$ws.Range("F37").Value = 75
$ws.Range("H37").Value = 112
$ws.Range("E41").Value = 387
$ws.Range("E42").Value = 379
$ws.Range("F42").Value = 208
$ws.Range("G42").Value = 61
$ws.Range("H42").Value = 269
$ws.Range("E43").Value = 118
$ws.Range("G43").Value = 27
$ws.Range("H43").Value = 90
$ws.Range("E44").Value = 311
$ws.Range("E46").Value = 322
$ws.Range("E47").Value = 453
$ws.Range("F47").Value = 225
$ws.Range("H47").Value = 317
$ws.Range("E48").Value = 212
$ws.Range("F48").Value = 91
$ws.Range("H48").Value = 135
$ws.Range("E51").Value = 233
